$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

Set-TextValue 'D2' '28.595.19'
Set-TextValue 'E2' '  +1.05%  '

Set-TextValue 'D3' '1.872.93'
Set-TextValue 'E3' '  +0.64%  '

Set-TextValue 'E4' '  -1.25%  '

Set-TextValue 'D5' '314.52'
Set-TextValue 'E5' '  -0.36%  '

Set-TextValue 'D6' '1.006'
Set-TextValue 'E6' '  -0.96%  '

Set-TextValue 'D7' '0.5082'
Set-TextValue 'E7' '  -0.20%  '

Set-TextValue 'D8' '0.3917'
Set-TextValue 'E8' '  -0.52%  '

Set-TextValue 'D9' '0.08351'
Set-TextValue 'E9' '  -0.99%  '

Set-TextValue 'E10' '  +1.24%  '

Set-TextValue 'D11' '1.105'
Set-TextValue 'E11' '  +0.04%  '

Set-TextValue 'D12' '6.192'
Set-TextValue 'E12' '  -0.59%  '

Set-TextValue 'D13' '1.865.86'
Set-TextValue 'E13' '  +3.45%  '

Set-TextValue 'D14' '20.33'
Set-TextValue 'E14' '  -0.05%  '

Set-TextValue 'E15' '  +1.16%  '

Set-TextValue 'D16' '1.007'
Set-TextValue 'E16' '  -1.12%  '

Set-TextValue 'D17' '93.11'
Set-TextValue 'E17' '  +3.10%  '

Set-TextValue 'E18' '  -0.93%  '

Set-TextValue 'D19' '0.06716'
Set-TextValue 'E19' '  -0.14%  '

Set-TextValue 'D20' '17.63'
Set-TextValue 'E20' '  +0.18%  '

Set-TextValue 'D21' '1.005'
Set-TextValue 'E21' '  -1.08%  '

Set-TextValue 'D22' '5.930'
Set-TextValue 'E22' '  +0.01%  '

Set-TextValue 'D23' '28.596.61'
Set-TextValue 'E23' '  +0.93%  '

Set-TextValue 'D24' '11.07'
Set-TextValue 'E24' '  -0.30%  '

Set-TextValue 'D25' '2.193'
Set-TextValue 'E25' '  -3.42%  '

Set-TextValue 'D26' '2.081.47'
Set-TextValue 'E26' '  +3.33%  '

Set-TextValue 'E27' '  -2.22%  '

Set-TextValue 'E28' '  -0.12%  '

Set-TextValue 'D29' '2.418'
Set-TextValue 'E29' '  +3.24%  '

Set-TextValue 'D30' '126.57'
Set-TextValue 'E30' '  -0.08%  '

Set-TextValue 'D31' '0.1036'
Set-TextValue 'E31' '  -1.03%  '

Set-TextValue 'D32' '1.044'
Set-TextValue 'E32' '  +1.49%  '

Set-TextValue 'D33' '5.773'
Set-TextValue 'E33' '  +0.48%  '

Set-TextValue 'D34' '3.632'
Set-TextValue 'E34' '  -0.17%  '

Set-TextValue 'D35' '0.02449'
Set-TextValue 'E35' '  +1.32%  '

Set-TextValue 'D36' '0.06546'
Set-TextValue 'E36' '  +1.83%  '

Set-TextValue 'D37' '9.032'
Set-TextValue 'E37' '  +2.73%  '

Set-TextValue 'D38' '0.2161'
Set-TextValue 'E38' '  -0.78%  '

Set-TextValue 'D39' '5.027'
Set-TextValue 'E39' '  +1.32%  '

Set-TextValue 'D40' '1.188'
Set-TextValue 'E40' '  +1.23%  '

Set-TextValue 'E41' '  -1.33%  '

Set-TextValue 'D42' '0.6367'
Set-TextValue 'E42' '  +0.35%  '

Set-TextValue 'E43' '  -0.41%  '

Set-TextValue 'E44' '  -0.71%  '

Set-TextValue 'D45' '0.5982'
Set-TextValue 'E45' '  -0.25%  '

Set-TextValue 'D46' '13.10'
Set-TextValue 'E46' '  +0.42%  '

Set-TextValue 'E47' '  -0.44%  '

Set-TextValue 'D48' '2.001'
Set-TextValue 'E48' '  +1.12%  '

Set-TextValue 'D49' '1.221'
Set-TextValue 'E49' '  +1.90%  '

Set-TextValue 'B50' 'WEMIXTOKEN'
Set-TextValue 'C50' 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
Set-TextValue 'D50' '1.197'
Set-TextValue 'E50' '  -1.00%  '

Set-TextValue 'B51' 'Quant'
Set-TextValue 'C51' 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
Set-TextValue 'D51' '122.21'
Set-TextValue 'E51' '  +1.41%  '
